$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 4 form-field rows:
#   row 1: Username / Text   / 1
#   row 2: Email    / Email  / 2
#   row 3: Confirm Email / Email / 3 / Eq(2)
#   row 4: Password / Password / 4
#
# We add a "Label" component row above Username, Email and Password (but
# NOT above Confirm Email), duplicating the field-name row and turning the
# duplicate's first copy into the new Label row. Working from the bottom
# row upward means every insertion point we still need to touch keeps its
# original row number.

# --- Password (row 4) ---------------------------------------------------
$ws.Rows("4:4").Copy()
$ws.Rows("5:5").Insert()
$ws.Range("B4").Value = "Label"
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 7

# --- Email (row 2) -------------------------------------------------------
$ws.Rows("2:2").Copy()
$ws.Rows("3:3").Insert()
$ws.Range("B2").Value = "Label"
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4

# --- Username (row 1) -----------------------------------------------------
$ws.Rows("1:1").Copy()
$ws.Rows("2:2").Insert()
$ws.Range("B1").Value = "Label"
$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 2

# --- Confirm Email row (now row 5): renumber + fix the Eq() formula label
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "Eq(4)"

# --- Cosmetic workbook/view metadata (selection moves along with edits) --
$ws.Range("J8").Select() | Out-Null
